# Weekly data refresh: prepend two new "Angeleno" price rows (most recent
# report dates) ahead of the existing "Ciruela" records, pushing every
# existing data row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 44 (the start of the
# data block whose rows all shift down by 2).
$ws.Rows("44:45").Insert()

# New row 44 - Angeleno / Especial
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 45044
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100103
$ws.Range("H44").Value = "Frutos de hueso (carozo)"
$ws.Range("I44").Value = 100103002
$ws.Range("J44").Value = "Ciruela"
$ws.Range("K44").Value = "Angeleno"
$ws.Range("L44").Value = "Especial"
$ws.Range("M44").Value = 50
$ws.Range("N44").Value = 12000
$ws.Range("O44").Value = 12000
$ws.Range("P44").Value = 12000
$ws.Range("Q44").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R44").Value = "Región de O'Higgins"
$ws.Range("S44").Value = 667
$ws.Range("T44").Value = 18

# New row 45 - Angeleno / Primera
$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 45044
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100103
$ws.Range("H45").Value = "Frutos de hueso (carozo)"
$ws.Range("I45").Value = 100103002
$ws.Range("J45").Value = "Ciruela"
$ws.Range("K45").Value = "Angeleno"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 50
$ws.Range("N45").Value = 10000
$ws.Range("O45").Value = 10000
$ws.Range("P45").Value = 10000
$ws.Range("Q45").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R45").Value = "Región de O'Higgins"
$ws.Range("S45").Value = 556
$ws.Range("T45").Value = 18

Write-Host "Inserted 2 rows with new Angeleno price data at rows 44-45; used range now has $($ws.UsedRange.Rows.Count) rows"
